$d = $word.ActiveDocument

# Update the date heading in the first paragraph
$dateRange = $d.Paragraphs.Item(1).Range
$dateRange.Find.Execute("2025-08-06 Wednesday", $true, $false, $false, $false, $false, $true, 0, $false, "2025-08-07 Thursday", 1) | Out-Null

# Update each table cell value in row-major order (old -> new pairs from the diff)
$t = $d.Tables.Item(1)
$values = @(
    @("81-62=", "59-42="),
    @("16+1=", "70-0="),
    @("55+38=", "44-22="),
    @("87+11=", "11+61="),
    @("86-83=", "41+22="),
    @("25+19=", "14+83="),
    @("92-81=", "76-21="),
    @("81-62=", "76-63="),
    @("14-4=", "45+31="),
    @("88-64=", "55+32="),
    @("16+16=", "81-66="),
    @("95-26=", "2+84="),
    @("27-27=", "24+25="),
    @("38-22=", "24+54="),
    @("57+23=", "31+49="),
    @("23+67=", "82-61="),
    @("82-55=", "47+40="),
    @("23-19=", "94-28="),
    @("83-8=", "94-28="),
    @("21-20=", "49+46="),
    @("80-65=", "92-92="),
    @("94-27=", "6+6="),
    @("80-2=", "87-30="),
    @("4+40=", "86-50="),
    @("26+41=", "69-64="),
    @("46+37=", "68-66="),
    @("0+11=", "17-1="),
    @("92-11=", "77-42="),
    @("68+20=", "79-55="),
    @("78-54=", "46+8="),
    @("54-36=", "78-3="),
    @("33+40=", "85+6="),
    @("0+43=", "28+44="),
    @("77+11=", "76-56="),
    @("25+71=", "61+18="),
    @("78-39=", "5+14="),
    @("73-12=", "99-12="),
    @("24+7=", "14+65="),
    @("66-52=", "13+45="),
    @("42+26=", "9+76="),
    @("77+18=", "4-4="),
    @("32+30=", "41+5="),
    @("20+2=", "37-12="),
    @("94-80=", "40+1="),
    @("95-6=", "53+24="),
    @("67-10=", "89-67="),
    @("32-7=", "73-5="),
    @("81-24=", "81-0="),
    @("19+25=", "69-28="),
    @("69+5=", "51+0="),
    @("64-10=", "82+0="),
    @("38+8=", "45-41="),
    @("67-43=", "36-15="),
    @("97-71=", "7+27="),
    @("5+29=", "34+22="),
    @("29+12=", "70-32="),
    @("31+18=", "68-57="),
    @("18+76=", "57+31="),
    @("74-48=", "47-46="),
    @("64-4=", "70-27="),
    @("58-55=", "42-30="),
    @("39+44=", "2+14="),
    @("39-8=", "84-3="),
    @("69-69=", "27-6="),
    @("46+37=", "40+33="),
    @("32-30=", "42+11="),
    @("22+27=", "17+24="),
    @("88-40=", "24-11="),
    @("36+16=", "52-1="),
    @("54+2=", "16+23="),
    @("60+19=", "75-47="),
    @("37+37=", "51+38="),
    @("29+41=", "27+47="),
    @("17+57=", "63+10="),
    @("9+1=", "87-7="),
    @("44-18=", "69-20="),
    @("27+42=", "15+33="),
    @("15+48=", "25-21="),
    @("40+15=", "73-4="),
    @("58+16=", "11-1="),
    @("30+44=", "51+37="),
    @("69-42=", "24+63="),
    @("45-5=", "8+31="),
    @("10-7=", "99-13="),
    @("35-27=", "53+2="),
    @("31+7=", "21+60="),
    @("1+66=", "77+7="),
    @("76-37=", "8+66="),
    @("8+50=", "88-9="),
    @("44-1=", "5+4="),
    @("18+48=", "53+34="),
    @("43+36=", "65+4="),
    @("71-27=", "83-0="),
    @("78-24=", "49+30="),
    @("82-34=", "89-57="),
    @("29+32=", "13-3="),
    @("79-39=", "78-23="),
    @("20+46=", "62-7="),
    @("20+12=", "71-34="),
    @("69-12=", "39-24=")
)

$numCols = 5
$index = 0
for ($row = 1; $row -le $t.Rows.Count; $row++) {
    for ($col = 1; $col -le $numCols; $col++) {
        $pair = $values[$index]
        $oldVal = $pair[0]
        $newVal = $pair[1]
        $cell = $t.Cell($row, $col)
        $cr = $cell.Range
        $cr.Find.Execute($oldVal, $true, $false, $false, $false, $false, $true, 0, $false, $newVal, 1) | Out-Null
        $index = $index + 1
    }
}

Write-Output "Done: updated $index cells"
